$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5150
$ws1.Range("F7").Value = 781
$ws1.Range("F8").Value = 266

# Sheet "演出" (shows)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 33

# Sheet "全部类型" (all types - combined)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5150
$ws4.Range("F7").Value = 781
$ws4.Range("F8").Value = 33
$ws4.Range("F9").Value = 266
